# Generate Report for Handback
# Adds a second handback row (839848ef-74c6-441f-b782-c0f541aa1a71) to the
# Overview / zh-cn / de-de sheets, alongside the existing
# (renamed) 6d63531b... -> 53478499... row.

$wb = $excel.ActiveWorkbook

$hyperlinkUnderline = 2          # xlUnderlineStyleSingle
$hyperlinkColor = 15570276       # BGR of FF6495ED (R=0x64,G=0x95,B=0xED)

function Set-HyperlinkStyle($range) {
    $range.Font.Underline = $hyperlinkUnderline
    $range.Font.Color = $hyperlinkColor
}

function Set-DateStyle($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

# The existing row's "Latest HO Xliff Generate Date" text is refreshed too.
$wsOverview.Range("G2").Value = "2016-09-05 13:16:12"

$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "839848ef-74c6-441f-b782-c0f541aa1a71.md"
$wsOverview.Range("B3").Value = "e2e\839848ef-74c6-441f-b782-c0f541aa1a71.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-09-05 13:16:12"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90a3c1c40685dab7066e7b2c01eccc4134bf4f91/e2e/839848ef-74c6-441f-b782-c0f541aa1a71.md", "", "", "e2e\839848ef-74c6-441f-b782-c0f541aa1a71.md") | Out-Null
Set-HyperlinkStyle $wsOverview.Range("B3")
Set-DateStyle $wsOverview.Range("G3")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)

# Refresh existing (renamed-file) row values.
$wsZhCn.Range("G2").Value = "53478499-934e-480e-bd4b-871f57b73d9c.e2160e0be45c77815f671bc1b7d30101bbb330bc.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-05 13:15:58"
$wsZhCn.Range("J2").Value = "53478499-934e-480e-bd4b-871f57b73d9c.e2160e0be45c77815f671bc1b7d30101bbb330bc.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-05 13:16:33"

$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = "839848ef-74c6-441f-b782-c0f541aa1a71.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'True"
$wsZhCn.Range("G3").Value = "839848ef-74c6-441f-b782-c0f541aa1a71.7432e153b21e08bcc6d0056b3d9978990a763004.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-05 13:15:58"
$wsZhCn.Range("I3").Value = "839848ef-74c6-441f-b782-c0f541aa1a71.md"
$wsZhCn.Range("J3").Value = "839848ef-74c6-441f-b782-c0f541aa1a71.7432e153b21e08bcc6d0056b3d9978990a763004.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-05 13:16:33"
$wsZhCn.Range("L3").Value = "'"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = "'"
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = "'"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90a3c1c40685dab7066e7b2c01eccc4134bf4f91/e2e/839848ef-74c6-441f-b782-c0f541aa1a71.md", "", "", "839848ef-74c6-441f-b782-c0f541aa1a71.md") | Out-Null
Set-HyperlinkStyle $wsZhCn.Range("A3")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ec95363a85cf68adcc4d3c7add5eb4925e30610c/e2e/839848ef-74c6-441f-b782-c0f541aa1a71.md", "", "", "839848ef-74c6-441f-b782-c0f541aa1a71.md") | Out-Null
Set-HyperlinkStyle $wsZhCn.Range("I3")

Set-DateStyle $wsZhCn.Range("H3")
Set-DateStyle $wsZhCn.Range("K3")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)

# Refresh existing (renamed-file) row values.
$wsDeDe.Range("G2").Value = "53478499-934e-480e-bd4b-871f57b73d9c.e2160e0be45c77815f671bc1b7d30101bbb330bc.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-05 13:16:12"
$wsDeDe.Range("J2").Value = "53478499-934e-480e-bd4b-871f57b73d9c.e2160e0be45c77815f671bc1b7d30101bbb330bc.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-05 13:16:41"

$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = "839848ef-74c6-441f-b782-c0f541aa1a71.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'True"
$wsDeDe.Range("G3").Value = "839848ef-74c6-441f-b782-c0f541aa1a71.7432e153b21e08bcc6d0056b3d9978990a763004.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-05 13:16:12"
$wsDeDe.Range("I3").Value = "839848ef-74c6-441f-b782-c0f541aa1a71.md"
$wsDeDe.Range("J3").Value = "839848ef-74c6-441f-b782-c0f541aa1a71.7432e153b21e08bcc6d0056b3d9978990a763004.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-05 13:16:41"
$wsDeDe.Range("L3").Value = "'"
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("N3").Value = "'"
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = "'"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90a3c1c40685dab7066e7b2c01eccc4134bf4f91/e2e/839848ef-74c6-441f-b782-c0f541aa1a71.md", "", "", "839848ef-74c6-441f-b782-c0f541aa1a71.md") | Out-Null
Set-HyperlinkStyle $wsDeDe.Range("A3")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/5b2ff27c96e28bf42bbfe643efb5994b56aafa7d/e2e/839848ef-74c6-441f-b782-c0f541aa1a71.md", "", "", "839848ef-74c6-441f-b782-c0f541aa1a71.md") | Out-Null
Set-HyperlinkStyle $wsDeDe.Range("I3")

Set-DateStyle $wsDeDe.Range("H3")
Set-DateStyle $wsDeDe.Range("K3")
